# Add a new "PasswordChangeComplete" validation row to the "Reset password"
# sheet: Label / Locator / ExpectedString for the post-reset confirmation
# message, so the test string-reader can locate it too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new row's values. Order matters here only insofar as it drives
# the order new entries land in the shared-string table (B then A then C).
$ws.Range("B2").Value = "//p[@class='ep-complete__info-text']"
$ws.Range("A2").Value = "PasswordChangeComplete"
$ws.Range("C2").Value = "パスワードの変更が完了しました。"

# Size the Label/Locator columns to fit their (now longer) contents, like
# the original author did by hand in Excel.
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
